# Apply the edits described by the diff to the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing columns I, K, M, O for data rows 2-25 ---
# Before: I=1, K=2, M=1, O=2  ->  After: I=2, K=1, M=2, O=1
for ($r = 2; $r -le 25; $r++) {
    $ws.Range("I$r").Value = 2
    $ws.Range("K$r").Value = 1
    $ws.Range("M$r").Value = 2
    $ws.Range("O$r").Value = 1
}

# --- Add new header cells P1 (14) and Q1 (15), matching style/format of existing header row ---
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Add new data columns P and Q for rows 2-25, all equal to 2 ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Range("P$r").Value = 2
    $ws.Range("Q$r").Value = 2
}
